$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F12").Value = -3
$ws.Range("F18").Value = -9
$ws.Range("F22").Value = -1
$ws.Range("F23").Value = -2
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = -3
$ws.Range("F36").Value = 0
